$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix D20: "excluido" -> "nao excluido"
$ws.Range("D20").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos com o Periodo Avaliativo nao excluido"

# 2) Fix the "Data Inicial e Data Final" wording in the repeated test-case blocks
#    (rows 50, 61, 72, 85 share the same two strings)
$rows = @(50, 61, 72, 85)
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "Lider de Pessoas preenche o campo 'Data Inicial' e 'Data Final' informando as respectivas datas referentes ao periodo"
    $ws.Range("D$r").Value = "SYSTEM apresenta o campo 'Data Inicial' e 'Data Final' preenchido corretamente"
}
